# Applies the symbol-list price/volume refresh described in the commit
# "Updated symbol list on Wed Jan  4 04:42:22 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (Price column D / Volume(1h) column E).
$updates = [ordered]@{
    "D2" = "256.06"
    "E2" = "4.27%"
    "D3" = "28.12"
    "E3" = "-4.33%"
    "D4" = "5.347"
    "E4" = "3.78%"
    "D5" = "0.05822"
    "E5" = "0.76%"
    "D6" = "6.708"
    "E6" = "1.43%"
    "D7" = "3.244"
    "E7" = "2.46%"
    "D8" = "0.8730"
    "E8" = "1.95%"
    "D9" = "0.9064"
    "E9" = "5.93%"
    "D10" = "0.1418"
    "E10" = "3.98%"
    "D11" = "0.07205"
    "E11" = "2.63%"
    "D12" = "0.03174"
    "E12" = "4.00%"
    "D13" = "0.09237"
    "E13" = "-1.47%"
    "D14" = "0.001549"
    "E14" = "1.17%"
    "D15" = "0.0006065"
    "E15" = "1.10%"
    "D16" = "0.006039"
    "E16" = "-0.09%"
    "E17" = "0.44%"
    "E18" = "4.63%"
    "E19" = "-1.15%"
    "D20" = "0.03470"
    "E20" = "4.50%"
    "D21" = "0.1310"
    "E21" = "2.20%"
    "D22" = "3.527"
    "E22" = "6.35%"
    "D23" = "0.04167"
    "E23" = "0.65%"
    "E24" = "-1.50%"
    "D25" = "0.001223"
    "E25" = "-0.22%"
    "D26" = "0.004881"
    "E26" = "18.08%"
    "E27" = "-0.87%"
    "E28" = "0.61%"
    "D40" = "0.03852"
    "D41" = "0.005729"
    "E41" = "-2.36%"
    "D42" = "0.1102"
    "E42" = "3.12%"
    "D43" = "0.002388"
    "E43" = "8.59%"
    "D44" = "0.009897"
    "E44" = "15.95%"
    "D45" = "0.00005260"
    "E45" = "-0.46%"
    "D46" = "0.00000000749"
    "E46" = "-0.06%"
    "D47" = "0.09991"
    "E47" = "72.28%"
    "D48" = "0.002136"
    "E48" = "-1.65%"
    "D49" = "0.00002098"
    "E49" = "-0.06%"
    "D50" = "0.0001998"
    "E50" = "-0.06%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text interpretation so values like "4.27%" or "0.8730" are stored
    # verbatim (same text they were before) instead of being parsed as numbers,
    # which would silently change formatting / drop trailing zeros.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Restore the default style so the cell looks just like it did originally
    # (no explicit number-format override left behind).
    $cell.Style = "Normal"
}
